$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 19: add "correct" test results in E19:G19 (style matches row 17's E:G pattern) ---
$ws.Range("E17:G17").Copy()
$ws.Range("E19:G19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E19").Value = "correct"
$ws.Range("F19").Value = "correct"

# --- Row 21: add "correct" test results in E21:F21 ---
$ws.Range("E17:F17").Copy()
$ws.Range("E21:F21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E21").Value = "correct"
$ws.Range("F21").Value = "correct"

# --- Row 22: add E22 ("correct", shaded style like E18/E20) and F22 (" correct", no special style) ---
$ws.Range("E18").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E22").Value = "correct"
$ws.Range("F22").Value = " correct"

# --- Row 23: D23 gets the shaded style (like E27/F27), E23:G23 follow the row17 pattern ---
$ws.Range("E27").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D23").Value = "yes"

$ws.Range("E17:G17").Copy()
$ws.Range("E23:G23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E23").Value = "correct"
$ws.Range("F23").Value = "correct"

# --- Row 24: E24 and F24 become "correct" with the shaded style used in row 36 ---
$ws.Range("E36:F36").Copy()
$ws.Range("E24:F24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E24").Value = "correct"
$ws.Range("F24").Value = "correct"

# --- Update the active selection to E22 ---
$ws.Range("E22").Select()
